$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Shape 1: content placeholder (內容版面配置區 2) ---
$sh1 = $s.Shapes.Item(1)

# Move the shape up (y offset 764704 -> 620688 EMU)
$sh1.Top = 48.8731

$tr1 = $sh1.TextFrame.TextRange

# Replace the whole text with a single paragraph first, so that the
# paragraph-level formatting (baseline alignment + custom bullet) is
# applied while there is still only one paragraph; splitting the text
# into multiple paragraphs afterwards copies this paragraph formatting
# to each new paragraph.
$tr1.Text = "首先神經網絡中需要有接收信號的主體，類似於樹突和樹突末梢。其次需要有一個類似於細胞體的結構來綜合處理接收到的信號。最後通過一個的機制就處理完成的信號發送出去，類似於軸突的工作。"

$pf1 = $tr1.ParagraphFormat
$pf1.BaseLineAlignment = 1
$bf1 = $pf1.Bullet
$bf1.Font.Name = "Wingdings"
$bf1.Type = 1
$bf1.Character = 108

# Now split the single paragraph into three separate paragraphs.
$tr1b = $sh1.TextFrame.TextRange
$tr1b.Text = "首先神經網絡中需要有接收信號的主體，類似於樹突和樹突末梢。`r其次需要有一個類似於細胞體的結構來綜合處理接收到的信號。`r最後通過一個的機制就處理完成的信號發送出去，類似於軸突的工作。"

# --- Shape 4: rectangle (矩形 5) ---
$sh4 = $s.Shapes.Item(4)
$tr4 = $sh4.TextFrame.TextRange

# Delete and re-insert the same text so the trailing empty-run
# (endParaRPr) marker left over from the previous edit is dropped.
$tr4.Delete()
$tr4b = $sh4.TextFrame.TextRange
$tr4b.Text = "此算法將可以使用非線性近似將資料分類或進行迴歸運算。"

# Re-typing the text triggers the shape's auto-fit height
# recalculation; restore the original height so the shape geometry is
# left untouched.
$sh4.Height = 50.8922
